# Updated cryptos list on Wed Mar 20 14:55:31 UTC 2024 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price values are plain numeric-looking strings (e.g. "0.610").
# Pre-format each target cell as Text so Excel keeps the exact literal
# string instead of silently coercing it to a Number (which would drop
# meaningful trailing zeros / introduce floating point noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.185.92"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "3.373.94"
$ws.Range("E3").Value = "  +3.37%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "525.19"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").Value = "174.92"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").Value = "3.367.48"
$ws.Range("E8").Value = "  +3.16%  "

$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").Value = "0.610"
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").Value = "53.50"
$ws.Range("E11").Value = "  -7.46%  "

$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").Value = "9.11"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "3.900.51"
$ws.Range("E15").Value = "  +3.41%  "

$ws.Range("D16").Value = "3.363.90"
$ws.Range("E16").Value = "  +3.39%  "

$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").Value = "17.61"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").Value = "64.023.97"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "11.30"
$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("D21").Value = "0.967"
$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").Value = "374.65"
$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("D23").Value = "11.61"
$ws.Range("E23").Value = "  +3.34%  "

$ws.Range("E24").Value = "  +8.07%  "

$ws.Range("D25").Value = "81.43"
$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").Value = "3.71"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").Value = "6.17"
$ws.Range("E27").Value = "  +1.52%  "

$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +3.09%  "

$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("D31").Value = "28.97"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("D32").Value = "631.67"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").Value = "6.47"
$ws.Range("E33").Value = "  -4.80%  "

$ws.Range("D34").Value = "11.24"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").Value = "58.12"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("D40").Value = "0.0₃0746"
$ws.Range("E40").Value = "  +14.07%  "

$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  +8.54%  "

$ws.Range("D43").Value = "2.976.95"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Value = "0.126"
$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("D45").Value = "2.98"
$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("E46").Value = "  +4.07%  "

$ws.Range("D47").Value = "0.0397"
$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("E49").Value = "  +2.66%  "

$ws.Range("D50").Value = "0.126"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "137.13"
$ws.Range("E51").Value = "  +4.69%  "
